# Fluid Capacity Mass.xlsx — apply author's edits:
#  - container_diameter (B4) value 10 -> 15
#  - fluid mass formula (B6) rewritten to use defined names, new cached value
#  - "fluid_mass" defined name broken to a #REF! (as in the authoritative diff)
#  - active selection moved from C8 to J11

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# container_diameter: 10 -> 15
$ws.Range("B4").Value = 15

# fluid mass formula now references the named ranges instead of raw cells
$ws.Range("B6").Formula = "=(container_height*3.14159*(container_diameter/2)^2)/1000"

# The "fluid_mass" defined name now points at a broken reference
$wb.Names.Item("fluid_mass").RefersTo = "=Sheet1!#REF!"

# Move the active selection to J11
$ws.Range("J11").Select()
